$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = -12.06379999999999
    "B10" = 6.1053
    "B12" = 5.496600000000002
    "C12" = -11.20629999999999
    "C17" = -13.32909999999999
    "B18" = 6.826299999999993
    "C26" = -12.51910000000001
    "C27" = -12.6265
    "C28" = -13.59819999999999
    "B37" = 8.609300000000006
    "C37" = -12.6472
    "B55" = 6.051399999999994
    "C65" = -12.8144
    "B68" = 4.797399999999996
    "C73" = -11.70700000000001
    "B77" = 9.480500000000005
    "B78" = 9.778999999999998
    "C84" = -13.26109999999999
    "C85" = -12.4904
    "C93" = -10.2789
    "C95" = -12.1492
    "C98" = -12.263
    "C99" = -11.85020000000001
    "C101" = -12.958
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
